$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns:
#  - one before the old "Study_Types" column (old col C)
#  - one before the old "Template_name" column (old col D, which after the
#    first insert has shifted to col E)
$ws.Columns("C:C").Insert()
$ws.Columns("E:E").Insert()

# New column C ("Population_Radio_button") - header + data, filled top to
# bottom first so the shared strings land in the right order.
$ws.Range("C1").Value = "Population_Radio_button"
$ws.Range("C2").Value = "Test - Test_radio_button"

# New column E ("slrtype_Radio_button") - header + data.
$ws.Range("E1").Value = "slrtype_Radio_button"
$ws.Range("E2").Value = "Clinical_radio_button"
$ws.Range("E3").Value = "Economic_radio_button"
$ws.Range("E4").Value = "Quality of Life_radio_button"
$ws.Range("E5").Value = "Real-world Evidence_radio_button"

# The new header cells should carry the default (unstyled) format, not the
# centered header style that Insert copied from their neighbour.
$ws.Range("C1").Style = "Normal"
$ws.Range("E1").Style = "Normal"

# Give the two new columns their own (narrower) widths, matching the
# neighbouring columns they sit next to.
$ws.Columns("C:C").ColumnWidth = 9
$ws.Columns("E:E").ColumnWidth = 17

# Update the active selection to match the saved workbook state
$ws.Range("K1").Select()
